# Updated Mobile automation script of Transfer note
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) org.open xpath index changed from 106 to 1
$ws.Range("B25").Value = "(//button[text()=' Open '])[1]"

# 2) SupplyDate value changed
$ws.Range("B97").Value = "22-02-2024"

# 3) receiptDate value changed (shared string also referenced by B131)
$ws.Range("B98").Value = "25-02-2024"
$ws.Range("B131").Value = "25-02-2024"

# 4) Append a new row 134 with the Mobile.data.xpath entry
$ws.Cells.Item(134, 1).Value = "Mobile.data.xpath"
$ws.Cells.Item(134, 2).Value = "//div[@class='modal-body px-0']/div/div[2]/div/div/div[2]/table/tbody/tr[7]/td[2]div/input"

# 5) Update active selection to the next empty row, matching the original author's cursor move
$ws.Range("A135").Select()
